# Work Diary update: append a new diary entry ("13/03/14" + note) after the
# last existing entry, keeping the hidden "_GoBack" bookmark pinned to the
# very end of the document (where Word leaves it after the last edit).

$d = $word.ActiveDocument

# 1. Grab the existing _GoBack bookmark (sits at the very end of the doc,
#    right after the last run of the last paragraph) and use its range to
#    insert the two new paragraphs *before* it.
$goBack = $d.Bookmarks.Item("_GoBack")
$insertionPoint = $goBack.Range
$insertionPoint.InsertBefore("`r13/03/14`rTried to encapsulate head into a List struct but ran into problems calling list->add in pthread_create")

# Inserting text in front of the bookmark re-seats it at the start of the
# newly typed final paragraph instead of the end, so it needs to be put
# back where Word would actually leave it: right after the last run of the
# (new) last paragraph.
$goBack2 = $d.Bookmarks.Item("_GoBack")
$goBack2.Delete()

$lastPara = $d.Paragraphs.Last

# 2. Temporarily add a scratch paragraph after the real last paragraph so
#    that the real last paragraph's end position becomes an "interior"
#    paragraph boundary. Adding a bookmark exactly at the end of a run
#    right before the final paragraph mark of the whole document doesn't
#    stick, but the same position *does* stick while more content follows.
$tailRange = $d.Content
$tailRange.Collapse(0)
$tailRange.InsertAfter("`rTEMP")

$boundary = $lastPara.Range.End
$bmRange = $d.Range($boundary, $boundary)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. Merge the scratch paragraph back into the real last paragraph (delete
#    its paragraph mark) -- the bookmark shifts left by one and now sits
#    exactly after the real text and before the scratch text.
$mergeRange = $d.Range($boundary - 1, $boundary)
$mergeRange.Delete()

# 4. Remove the scratch "TEMP" text, leaving the bookmark in place right
#    after the real last run.
$scratchRange = $d.Range($boundary - 1, $boundary - 1 + 4)
$scratchRange.Delete()
